# Update IFRS financial figures for 현대백화점 (rows 2-9) to corrected values.
# Commit message: "error solve ifrs list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 15519
$ws.Range("E2").Value = 3637
$ws.Range("F2").Value = 3637
$ws.Range("G2").Value = 3964
$ws.Range("H2").Value = 2910
$ws.Range("I2").Value = 2684
$ws.Range("J2").Value = 226
$ws.Range("K2").Value = 56839
$ws.Range("L2").Value = 18949
$ws.Range("M2").Value = 37890
$ws.Range("N2").Value = 32443
$ws.Range("O2").Value = 5446
$ws.Range("P2").Value = 1170
$ws.Range("Q2").Value = 3985
$ws.Range("R2").Value = -1934
$ws.Range("S2").Value = -2670
$ws.Range("T2").Value = 3717
$ws.Range("U2").Value = 268
$ws.Range("V2").Value = 5006
$ws.Range("W2").Value = 23.43
$ws.Range("X2").Value = 18.75
$ws.Range("Y2").Value = 8.57
$ws.Range("Z2").Value = 5.16
$ws.Range("AA2").Value = 50.01
$ws.Range("AB2").Value = 2699.91
$ws.Range("AC2").Value = 11470
$ws.Range("AD2").Value = 10.72
$ws.Range("AE2").Value = 140843
$ws.Range("AF2").Value = 0.87
$ws.Range("AG2").Value = 700
$ws.Range("AH2").Value = 0.57
$ws.Range("AI2").Value = 6.01
$ws.Range("AJ2").Value = 23402441
$ws.Range("D3").Value = 16570
$ws.Range("E3").Value = 3628
$ws.Range("F3").Value = 3628
$ws.Range("G3").Value = 3887
$ws.Range("H3").Value = 2803
$ws.Range("I3").Value = 2409
$ws.Range("J3").Value = 394
$ws.Range("K3").Value = 61608
$ws.Range("L3").Value = 21300
$ws.Range("M3").Value = 40309
$ws.Range("N3").Value = 34499
$ws.Range("O3").Value = 5810
$ws.Range("P3").Value = 1170
$ws.Range("Q3").Value = 4011
$ws.Range("R3").Value = -5421
$ws.Range("S3").Value = 1673
$ws.Range("T3").Value = 5012
$ws.Range("U3").Value = -1001
$ws.Range("V3").Value = 7195
$ws.Range("W3").Value = 21.9
$ws.Range("X3").Value = 16.91
$ws.Range("Y3").Value = 7.2
$ws.Range("Z3").Value = 4.73
$ws.Range("AA3").Value = 52.84
$ws.Range("AB3").Value = 2889.96
$ws.Range("AC3").Value = 10294
$ws.Range("AD3").Value = 12.29
$ws.Range("AE3").Value = 151369
$ws.Range("AF3").Value = 0.84
$ws.Range("AG3").Value = 700
$ws.Range("AH3").Value = 0.55
$ws.Range("AI3").Value = 6.62
$ws.Range("AJ3").Value = 23402441
$ws.Range("D4").Value = 18318
$ws.Range("E4").Value = 3832
$ws.Range("F4").Value = 3832
$ws.Range("G4").Value = 4359
$ws.Range("H4").Value = 3211
$ws.Range("I4").Value = 2758
$ws.Range("J4").Value = 453
$ws.Range("K4").Value = 65873
$ws.Range("L4").Value = 22767
$ws.Range("M4").Value = 43106
$ws.Range("N4").Value = 36938
$ws.Range("O4").Value = 6168
$ws.Range("P4").Value = 1170
$ws.Range("Q4").Value = 4420
$ws.Range("R4").Value = -5244
$ws.Range("S4").Value = 1084
$ws.Range("T4").Value = 3335
$ws.Range("U4").Value = 1085
$ws.Range("V4").Value = 8591
$ws.Range("W4").Value = 20.92
$ws.Range("X4").Value = 17.53
$ws.Range("Y4").Value = 7.72
$ws.Range("Z4").Value = 5.04
$ws.Range("AA4").Value = 52.82
$ws.Range("AB4").Value = 3116.85
$ws.Range("AC4").Value = 11784
$ws.Range("AD4").Value = 9.25
$ws.Range("AE4").Value = 162716
$ws.Range("AF4").Value = 0.67
$ws.Range("AG4").Value = 700
$ws.Range("AH4").Value = 0.64
$ws.Range("AI4").Value = 5.77
$ws.Range("AJ4").Value = 23402441
$ws.Range("D5").Value = 18481
$ws.Range("E5").Value = 3937
$ws.Range("F5").Value = 3937
$ws.Range("G5").Value = 4456
$ws.Range("H5").Value = 3022
$ws.Range("I5").Value = 2537
$ws.Range("J5").Value = 485
$ws.Range("K5").Value = 67773
$ws.Range("L5").Value = 21596
$ws.Range("M5").Value = 46177
$ws.Range("N5").Value = 39554
$ws.Range("O5").Value = 6623
$ws.Range("P5").Value = 1170
$ws.Range("Q5").Value = 4372
$ws.Range("R5").Value = -1457
$ws.Range("S5").Value = -2751
$ws.Range("T5").Value = 3144
$ws.Range("U5").Value = 1228
$ws.Range("V5").Value = 6197
$ws.Range("W5").Value = 21.3
$ws.Range("X5").Value = 16.35
$ws.Range("Y5").Value = 6.63
$ws.Range("Z5").Value = 4.52
$ws.Range("AA5").Value = 46.77
$ws.Range("AB5").Value = 3332.35
$ws.Range("AC5").Value = 10841
$ws.Range("AD5").Value = 9.640000000000001
$ws.Range("AE5").Value = 175349
$ws.Range("AF5").Value = 0.6
$ws.Range("AG5").Value = 800
$ws.Range("AH5").Value = 0.77
$ws.Range("AI5").Value = 7.11
$ws.Range("AJ5").Value = 23402441
$ws.Range("D6").Value = 18622
$ws.Range("E6").Value = 3567
$ws.Range("F6").Value = 3567
$ws.Range("G6").Value = 3983
$ws.Range("H6").Value = 2874
$ws.Range("I6").Value = 2390
$ws.Range("K6").Value = 69886
$ws.Range("L6").Value = 21595
$ws.Range("M6").Value = 48291
$ws.Range("N6").Value = 41325
$ws.Range("P6").Value = 1170
$ws.Range("Q6").Value = 4353
$ws.Range("R6").Value = -1341
$ws.Range("S6").Value = -1465
$ws.Range("T6").Value = 2213
$ws.Range("U6").Value = 2140
$ws.Range("V6").Value = 4991
$ws.Range("W6").Value = 19.15
$ws.Range("X6").Value = 15.43
$ws.Range("Y6").Value = 5.91
$ws.Range("Z6").Value = 4.18
$ws.Range("AA6").Value = 44.72
$ws.Range("AB6").Value = 3515.01
$ws.Range("AC6").Value = 10211
$ws.Range("AD6").Value = 8.85
$ws.Range("AE6").Value = 183201
$ws.Range("AF6").Value = 0.49
$ws.Range("AG6").Value = 900
$ws.Range("AH6").Value = 1
$ws.Range("AI6").Value = 8.5
$ws.Range("AJ6").Value = 23402441
$ws.Range("D7").Value = 22142
$ws.Range("E7").Value = 2880
$ws.Range("G7").Value = 3498
$ws.Range("H7").Value = 2518
$ws.Range("I7").Value = 2035
$ws.Range("K7").Value = 76500
$ws.Range("L7").Value = 26245
$ws.Range("M7").Value = 50256
$ws.Range("N7").Value = 42940
$ws.Range("P7").Value = 1170
$ws.Range("Q7").Value = 4425
$ws.Range("R7").Value = -5097
$ws.Range("S7").Value = 578
$ws.Range("T7").Value = 3379
$ws.Range("U7").Value = 1731
$ws.Range("W7").Value = 13.01
$ws.Range("X7").Value = 11.37
$ws.Range("Y7").Value = 4.83
$ws.Range("Z7").Value = 3.44
$ws.Range("AA7").Value = 52.22
$ws.Range("AC7").Value = 8698
$ws.Range("AD7").Value = 8.970000000000001
$ws.Range("AE7").Value = 192354
$ws.Range("AF7").Value = 0.41
$ws.Range("AG7").Value = 938
$ws.Range("AH7").Value = 1.2
$ws.Range("AI7").Value = 10.78
$ws.Range("D8").Value = 26124
$ws.Range("E8").Value = 3301
$ws.Range("G8").Value = 3871
$ws.Range("H8").Value = 2829
$ws.Range("I8").Value = 2360
$ws.Range("K8").Value = 80905
$ws.Range("L8").Value = 28195
$ws.Range("M8").Value = 52710
$ws.Range("N8").Value = 45047
$ws.Range("P8").Value = 1170
$ws.Range("Q8").Value = 4437
$ws.Range("R8").Value = -4066
$ws.Range("S8").Value = 166
$ws.Range("T8").Value = 3469
$ws.Range("U8").Value = 805
$ws.Range("W8").Value = 12.63
$ws.Range("X8").Value = 10.83
$ws.Range("Y8").Value = 5.36
$ws.Range("Z8").Value = 3.59
$ws.Range("AA8").Value = 53.49
$ws.Range("AC8").Value = 10083
$ws.Range("AD8").Value = 7.74
$ws.Range("AE8").Value = 201795
$ws.Range("AF8").Value = 0.39
$ws.Range("AG8").Value = 975
$ws.Range("AH8").Value = 1.25
$ws.Range("AI8").Value = 9.67
$ws.Range("D9").Value = 29156
$ws.Range("E9").Value = 3960
$ws.Range("G9").Value = 4545
$ws.Range("H9").Value = 3328
$ws.Range("I9").Value = 2786
$ws.Range("K9").Value = 84653
$ws.Range("L9").Value = 29107
$ws.Range("M9").Value = 55546
$ws.Range("N9").Value = 47527
$ws.Range("P9").Value = 1170
$ws.Range("Q9").Value = 5036
$ws.Range("R9").Value = -3630
$ws.Range("S9").Value = -665
$ws.Range("T9").Value = 3030
$ws.Range("U9").Value = 1710
$ws.Range("W9").Value = 13.58
$ws.Range("X9").Value = 11.42
$ws.Range("Y9").Value = 6.02
$ws.Range("Z9").Value = 4.02
$ws.Range("AA9").Value = 52.4
$ws.Range("AC9").Value = 11906
$ws.Range("AD9").Value = 6.55
$ws.Range("AE9").Value = 212904
$ws.Range("AF9").Value = 0.37
$ws.Range("AG9").Value = 1031
$ws.Range("AH9").Value = 8.66
$ws.Range("AI9").Value = 8.66
